$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (moved one month forward: 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# "COMUN" price table (rows 35-39)
$ws.Range("D35").Value = 3666.444
$ws.Range("D36").Value = 4110
$ws.Range("D37").Value = 4110
$ws.Range("D38").Value = 4110
$ws.Range("D39").Value = 10440.086

# "CON TOPE" price table (rows 42-46)
$ws.Range("D42").Value = 3953.09
$ws.Range("D43").Value = 4350
$ws.Range("D44").Value = 4350
$ws.Range("D45").Value = 4350

# D46's style carries a stale "quote prefix" flag (xf 13) left over from
# before it held a numeric value. Writing a new number through .Value
# makes Excel drop that flag (xf 9), same as typing a number over a
# quote-prefixed cell does interactively. Preserve the original
# formatting by stashing it on a scratch cell, updating the value, then
# pasting the formatting back.
$scratch = $ws.Range("F50")
$ws.Range("D46").Copy()
$scratch.PasteSpecial(-4122)
$ws.Range("D46").Value = 11676.446
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4122)
$scratch.Clear()
